# Module4Output_20240104.xlsx - ProductionPlan sheet update
# fix(module3): use uncon_planned_qty for future production; keep produced for today
#   - update MAT_A row quantities (future production now uses uncon_planned_qty)
#   - add a new MAT_B / LINE_B production row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductionPlan")

# --- Update existing row 2 (MAT_A) quantities ---
$ws.Range("G2").Value = 710
$ws.Range("H2").Value = 710
$ws.Range("J2").Value = 675

# --- Append new row 3 (MAT_B) ---
$ws.Range("A3").Value = "MAT_B"
$ws.Range("B3").Value = "PLANT_001"
$ws.Range("C3").Value = "LINE_B"
$ws.Range("D3").Value = 45295
$ws.Range("E3").Value = 45296
$ws.Range("F3").Value = 45297
$ws.Range("G3").Value = 104
$ws.Range("H3").Value = 104
$ws.Range("J3").Value = 92

# Apply the same date formatting used in row 2 to the new row's date cells
$ws.Range("D2:F2").Copy()
$ws.Range("D3:F3").PasteSpecial(-4122)
